# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) Slide Master, currently the
#                             "Integral" colour scheme.
#   ppt/theme/theme2.xml  -> bound only to the Notes Master, currently the
#                             default "Office Theme" colour scheme.
# The target revision swaps the two themes' content: the Slide Master theme
# becomes the stock "Office Theme" palette, and the Notes Master theme
# becomes the "Integral" palette.
#
# The PowerPoint object model only exposes a single live Theme (reachable
# from $p.SlideMaster.Theme / $p.NotesMaster.Theme / $p.DocumentTheme, which
# all resolve to the same underlying theme), so we drive the swap through
# its ThemeColorScheme, which is the editable surface for a theme's colours.

$p = $ppt.ActivePresentation

$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB() integers (0x00BBGGRR), replacing the current
# Integral palette so the Slide Master theme matches the post-edit
# "Office Theme" colours.
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
